$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Update the PRESUPUESTO column (G) for each advisor/client row.
$ws.Range("G3").Value = 500
$ws.Range("G5").Value = 3000
$ws.Range("G6").Value = 4000
$ws.Range("G7").Value = 1000
$ws.Range("G8").Value = 4000
$ws.Range("G11").Value = 7000
$ws.Range("G13").Value = 5000
$ws.Range("G14").Value = 7000
$ws.Range("G15").Value = 500
$ws.Range("G16").Value = 500
$ws.Range("G17").Value = 9500
$ws.Range("G20").Value = 2000
$ws.Range("G21").Value = 2500
$ws.Range("G22").Value = 6000

# Update the total row to reflect the new sum of the column above.
$ws.Range("G23").Value = 52500
